$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.907.55"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = "'2.357.89"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'0.693"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.95%  '
$ws.Range('D6').Value = "'242.36"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.30%  '
$ws.Range('D7').Value = "'76.64"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +4.71%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').Value = "'0.630"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +20.46%  '
$ws.Range('E10').Value = '  +4.43%  '
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('D12').Value = "'33.62"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +23.46%  '
$ws.Range('E13').Value = '  +14.25%  '
$ws.Range('E14').Value = '  +1.73%  '
$ws.Range('D15').Value = "'2.708.88"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').Value = "'16.89"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.58%  '
$ws.Range('D17').Value = "'0.935"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +5.61%  '
$ws.Range('D18').Value = "'2.362.00"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').Value = "'43.851.34"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('E20').Value = '  +2.54%  '
$ws.Range('D21').Value = "'6.70"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.46%  '
$ws.Range('D22').Value = "'77.80"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.97%  '
$ws.Range('D23').Value = "'263.11"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +4.78%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = "'2.53"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.08%  '
$ws.Range('D26').Value = "'3.63"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.83%  '
$ws.Range('E27').Value = '  +8.27%  '
$ws.Range('D28').Value = "'1.80"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +17.80%  '
$ws.Range('E29').Value = '  +2.45%  '
$ws.Range('D30').Value = "'23.13"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.03%  '
$ws.Range('D31').Value = "'175.24"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('E32').Value = '  -3.83%  '
$ws.Range('E33').Value = '  +4.63%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = "'5.42"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.79%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = "'0.0766"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +9.23%  '
$ws.Range('D36').Value = "'5.45"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +6.40%  '
$ws.Range('D37').Value = "'3.82"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.87%  '
$ws.Range('D38').Value = "'2.45"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.58%  '
$ws.Range('D39').Value = "'6.45"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.15%  '
$ws.Range('E40').Value = '  +7.33%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = "'0.211"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +19.39%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').Value = "'19.36"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = "'9.18"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.55%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = "'0.107"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +10.09%  '
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('D46').Value = "'2.54"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +11.76%  '
$ws.Range('D47').Value = "'1.26"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +4.25%  '
$ws.Range('E48').Value = '  +2.20%  '
$ws.Range('E49').Value = '  +1.91%  '
$ws.Range('D50').Value = "'4.58"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.00%  '
$ws.Range('D51').Value = "'56.30"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +10.55%  '
